$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.226.99"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.50"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.18"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06548"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.79"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07929"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.77"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.00"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.179"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6823"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.40"
$ws.Range("E16").Value = "  -5.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.229.83"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.77"
$ws.Range("E18").Value = "  +8.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007436"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.115.31"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.17"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.230"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.962"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.386"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.386"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.076"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04712"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7033"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.619"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.251"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.71"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.943"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8462"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4177"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.28"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "955.74"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.176"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.257"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.18"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +0.56%  "
